# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 289-290) into the Repollo sheet,
# pushing the existing rows 289-308 down to 291-310.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 289.
$ws.Rows("289:290").Insert()

# --- New row 289 -----------------------------------------------------
$ws.Cells.Item(289, 1).Value = 7
$ws.Cells.Item(289, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(289, 3).Value = "Ñuble"
$ws.Cells.Item(289, 4).Value = 44931
$ws.Cells.Item(289, 5).Value = 16
$ws.Cells.Item(289, 6).Value = 100112006
$ws.Cells.Item(289, 7).Value = "Repollo"
$ws.Cells.Item(289, 8).Value = "Crespo record"
$ws.Cells.Item(289, 9).Value = "Primera"
$ws.Cells.Item(289, 10).Value = 500
$ws.Cells.Item(289, 11).Value = 1200
$ws.Cells.Item(289, 12).Value = 1300
$ws.Cells.Item(289, 13).Value = 1250
$ws.Cells.Item(289, 14).Value = "$/unidad"
$ws.Cells.Item(289, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(289, 16).Value = 1250
$ws.Cells.Item(289, 17).Value = 1
$ws.Cells.Item(289, 18).Value = "Hortaliza"

# --- New row 290 -----------------------------------------------------
$ws.Cells.Item(290, 1).Value = 7
$ws.Cells.Item(290, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(290, 3).Value = "Ñuble"
$ws.Cells.Item(290, 4).Value = 44931
$ws.Cells.Item(290, 5).Value = 16
$ws.Cells.Item(290, 6).Value = 100112006
$ws.Cells.Item(290, 7).Value = "Repollo"
$ws.Cells.Item(290, 8).Value = "Crespo record"
$ws.Cells.Item(290, 9).Value = "Segunda"
$ws.Cells.Item(290, 10).Value = 600
$ws.Cells.Item(290, 11).Value = 1000
$ws.Cells.Item(290, 12).Value = 1100
$ws.Cells.Item(290, 13).Value = 1050
$ws.Cells.Item(290, 14).Value = "$/unidad"
$ws.Cells.Item(290, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(290, 16).Value = 1050
$ws.Cells.Item(290, 17).Value = 1
$ws.Cells.Item(290, 18).Value = "Hortaliza"
